# "document mis a jour" - add 4 new "Journal de travail" entries (rows 35-38)
# describing the "Calcul des points" / "Qui commence" / "Atout" features and
# the related Game controller update, on top of the existing log rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing, similarly-formatted row (row 22: date/activite/duree/etc.
# with the "React native (Frontend)" look - same fills/fonts/wrap as the new
# rows use) as the template for the formatting of the 4 new rows, then
# overwrite the values/formulas with the new journal entries.
$ws.Range("A22:H22").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Range("A36:H36").PasteSpecial(-4122)
$ws.Range("A37:H37").PasteSpecial(-4122)
$ws.Range("A38:H38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 35 - "Feature - Calcule des points"
$ws.Range("A35").Value = 44263
$ws.Range("B35").Value = "React native (Frontend)"
$ws.Range("C35").Value = "2h"
$ws.Range("D35").Value = "Feature - Calcule des points"
$ws.Range("E35").Value = "Oui"
$ws.Range("F35").Value = "Permet de calculer les points de chaque equipe et de renvoyer les points sur le serveur web"
$ws.Range("G35").Value = "Non"
$ws.Range("H35").Value = ""

# Row 36 - "Feature - Qui commence ?"
$ws.Range("A36").Value = 44263
$ws.Range("B36").Value = "React native (Frontend)"
$ws.Range("C36").Value = "2h"
$ws.Range("D36").Value = "Feature - Qui commence ?"
$ws.Range("E36").Value = "Oui"
$ws.Range("F36").Value = "Permet de savoir qui doit commencer et choisir l'atout"
$ws.Range("G36").Value = "Non"
$ws.Range("H36").Value = ""

# Row 37 - "Feature - Atout de la partie ?"
$ws.Range("A37").Value = 44263
$ws.Range("B37").Value = "React native (Frontend)"
$ws.Range("C37").Value = "1h"
$ws.Range("D37").Value = "Feature - Atout de la partie ?"
$ws.Range("E37").Value = "Oui"
$ws.Range("F37").Value = "Permet de savoir quel atout est présent pour la manche"
$ws.Range("G37").Value = "Non"
$ws.Range("H37").Value = ""

# Row 38 - "Mise a jour controller Game"
$ws.Range("A38").Value = 44263
$ws.Range("B38").Value = "React native (Frontend)"
$ws.Range("C38").Value = "1h"
$ws.Range("D38").Value = "Mise a jour controller Game"
$ws.Range("E38").Value = "Oui"
$ws.Range("F38").Value = "Mise  a jour du controller pour les manche et l'atout d'une partie"
$ws.Range("G38").Value = "Non"
$ws.Range("H38").Value = ""

# New rows render on two lines given the column widths/wrap - match the
# author's row heights.
$ws.Rows.Item(35).RowHeight = 31.5
$ws.Rows.Item(36).RowHeight = 31.5
$ws.Rows.Item(37).RowHeight = 31.5
$ws.Rows.Item(38).RowHeight = 31.5

# Keep the date column formatted the same way as the rest of the sheet
# (dd.mm.yy) rather than a raw serial number.
$ws.Range("A35:A38").NumberFormat = "dd\.mm\.yy;@"

# Scroll the sheet roughly to where the new rows were added, matching the
# author's final view state.
$ws.Range("G44").Select()
